$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-13 (A-D columns are text labels, E-T are numeric values)
$data = @(
    @("ECs",  "Spon2", "Itga5", "ECs",  2, 0.6666666666666666, 0.436226, 1.308678, 0.02612337021374315, 0.02612337021374314, 3, 1, 30.99161333333333, 92.97484, 0.3599121977633812, 0.3599121977633811, 13.51934751794667, 121.67412766152, 0.009402119586614747, 0.009402119586614743),
    @("ECs",  "Spon2", "Itga5", "FAPs", 2, 0.6666666666666666, 0.436226, 1.308678, 0.02612337021374315, 0.02612337021374314, 3, 1, 29.913269, 89.739807, 0.3473891556493311, 0.3473891556493311, 13.048945682794, 117.440511145146, 0.009074975521267119, 0.009074975521267117),
    @("ECs",  "Spon2", "Itga5", "M2",   2, 0.6666666666666666, 0.436226, 1.308678, 0.02612337021374315, 0.02612337021374314, 3, 1, 19.150218, 57.450654, 0.2223955550134164, 0.2223955550134163, 8.353822997267999, 75.184406975412, 0.005809721417506358, 0.005809721417506355),
    @("ECs",  "Spon2", "Itga5", "sCs",  2, 0.6666666666666666, 0.436226, 1.308678, 0.02612337021374315, 0.02612337021374314, 3, 1, 6.053716000000001, 18.161148, 0.07030309157387134, 0.07030309157387132, 2.640788315816, 23.767094842344, 0.001836553688354927, 0.001836553688354927),
    @("FAPs", "Spon2", "Itga5", "ECs",  3, 1, 15.76143266666667, 47.284298, 0.9438725354525366, 0.9438725354525365, 3, 1, 30.99161333333333, 92.97484, 0.3599121977633812, 0.3599121977633811, 488.4722267847022, 4396.25004106232, 0.3397112386432174, 0.3397112386432173),
    @("FAPs", "Spon2", "Itga5", "FAPs", 3, 1, 15.76143266666667, 47.284298, 0.9438725354525366, 0.9438725354525365, 3, 1, 29.913269, 89.739807, 0.3473891556493311, 0.3473891556493311, 471.4759751833873, 4243.283776650486, 0.32789108313145, 0.32789108313145),
    @("FAPs", "Spon2", "Itga5", "M2",   3, 1, 15.76143266666667, 47.284298, 0.9438725354525366, 0.9438725354525365, 3, 1, 19.150218, 57.450654, 0.2223955550134164, 0.2223955550134163, 301.834871558988, 2716.513844030892, 0.2099130563838874, 0.2099130563838874),
    @("FAPs", "Spon2", "Itga5", "sCs",  3, 1, 15.76143266666667, 47.284298, 0.9438725354525366, 0.9438725354525365, 3, 1, 6.053716000000001, 18.161148, 0.07030309157387134, 0.07030309157387132, 95.41523711712267, 858.737134054104, 0.0663571572939818, 0.06635715729398178),
    @("sCs",  "Spon2", "Itga5", "ECs",  3, 1, 0.5010290000000001, 1.503087, 0.0300040943337204, 0.0300040943337204, 3, 1, 30.99161333333333, 92.97484, 0.3599121977633812, 0.3599121977633811, 15.52769703678667, 139.74927333108, 0.01079883953354912, 0.01079883953354912),
    @("sCs",  "Spon2", "Itga5", "FAPs", 3, 1, 0.5010290000000001, 1.503087, 0.0300040943337204, 0.0300040943337204, 3, 1, 29.913269, 89.739807, 0.3473891556493311, 0.3473891556493311, 14.987415253801, 134.886737284209, 0.01042309699661401, 0.01042309699661401),
    @("sCs",  "Spon2", "Itga5", "M2",   3, 1, 0.5010290000000001, 1.503087, 0.0300040943337204, 0.0300040943337204, 3, 1, 19.150218, 57.450654, 0.2223955550134164, 0.2223955550134163, 9.594814574322001, 86.35333116889801, 0.006672777212022652, 0.006672777212022648),
    @("sCs",  "Spon2", "Itga5", "sCs",  3, 1, 0.5010290000000001, 1.503087, 0.0300040943337204, 0.0300040943337204, 3, 1, 6.053716000000001, 18.161148, 0.07030309157387134, 0.07030309157387132, 3.033087273764, 27.297785463876, 0.00210938059153462, 0.002109380591534619)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $rowNum = $startRow + $i
    $rowValues = $data[$i]
    for ($c = 0; $c -lt $rowValues.Count; $c++) {
        $ws.Cells.Item($rowNum, $c + 1).Value = $rowValues[$c]
    }
}
